$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 281
$ws1.Range("F5").Value = 3410
$ws1.Range("F8").Value = 162
$ws1.Range("F9").Value = 49
$ws1.Range("F13").Value = 1552

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 281
$ws4.Range("F5").Value = 3410
$ws4.Range("F9").Value = 162
$ws4.Range("F10").Value = 49
$ws4.Range("F16").Value = 1552
